$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("This report will also act as the readme file...GitHub"):
#    - merge the "...on my " / "GitHub" runs (drops the stray proofErr marks)
#    - append ". You can find the code for this homework listed "
#    - add a "here" hyperlink to the GitHub repo right after that
# ---------------------------------------------------------------------------

$rngMerge = $d.Content
$null = $rngMerge.Find.Execute(
    "This report will also act as the readme file that is posted on my GitHub",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This report will also act as the readme file that is posted on my GitHub", 2)

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2end = $d.Range($r2.Start, $r2.End - 1)
$r2end.Collapse(0)
$r2end.InsertAfter(". You can find the code for this homework listed ")
$r2end.Collapse(0)
$r2end.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$r3start = $d.Range($r3.Start, $r3.Start)
$null = $d.Hyperlinks.Add($r3start, "https://github.com/chomayouni/NLP/tree/main/Homework_1", "", "", "here")

$p2b = $d.Paragraphs.Item(2)
$r2b = $p2b.Range
$mark1 = $d.Range($r2b.End - 1, $r2b.End)
$mark1.Delete()

# ---------------------------------------------------------------------------
# 2) Remove the stray empty paragraph right after "...sections below."
# ---------------------------------------------------------------------------

$rngFind = $d.Content
$ok = $rngFind.Find.Execute("I" + [char]0x2019 + "ll talk about this in the sections below.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $hostPara = $rngFind.Paragraphs.Item(1)
    $emptyPara = $hostPara.Next()
    if ($emptyPara.Range.Text -eq [char]0x0D) {
        $emptyPara.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 3) Append a new paragraph at the end with a backup repository link, plus a
#    trailing empty paragraph.
# ---------------------------------------------------------------------------

$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $lastP.Range
$rLast.Collapse(0)
$rLast.InsertParagraphAfter()

$backupIdx = $d.Paragraphs.Count
$backupP = $d.Paragraphs.Item($backupIdx)
$backupP.Style = "Normal"
$backupP.Range.ListFormat.RemoveNumbers()
$backupP.Range.InsertAfter("Here is a backup repository link just in case: ")

$backupP2 = $d.Paragraphs.Item($backupIdx)
$rBack = $backupP2.Range
$rBackEnd = $d.Range($rBack.Start, $rBack.End - 1)
$rBackEnd.Collapse(0)
$rBackEnd.InsertParagraphAfter()

$linkHostIdx = $d.Paragraphs.Count
$linkHost = $d.Paragraphs.Item($linkHostIdx)
$linkHost.Style = "Normal"
$linkHost.Range.ListFormat.RemoveNumbers()
$rLinkHost = $linkHost.Range
$rLinkHostStart = $d.Range($rLinkHost.Start, $rLinkHost.Start)
$null = $d.Hyperlinks.Add($rLinkHostStart, "https://github.com/chomayouni/NLP/tree/main/Homework_1", "", "", "https://github.com/chomayouni/NLP/tree/main/Homework_1")

$backupP3 = $d.Paragraphs.Item($backupIdx)
$rBack3 = $backupP3.Range
$mark2 = $d.Range($rBack3.End - 1, $rBack3.End)
$mark2.Delete()

$veryLastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$rVeryLast = $veryLastP.Range
$rVeryLast.Collapse(0)
$rVeryLast.InsertParagraphAfter()
$finalEmpty = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalEmpty.Style = "Normal"
$finalEmpty.Range.ListFormat.RemoveNumbers()

Write-Host "Done."
